$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report week dates) ---
$volCell = $ws.Cells.Item(8, 1)
$volCell.Characters(21, 2).Text = "45"

$weekCell = $ws.Cells.Item(9, 3)
$weekCell.Characters(48, 9).Text = "11/12/2023"
$weekCell.Characters(27, 10).Text = "11/6/2023"

# --- Simple numeric value updates (style/format unchanged) ---
$ws.Cells.Item(15, 6).Value = 2
$ws.Cells.Item(15, 8).Value = 100
$ws.Cells.Item(15, 10).Value = 23
$ws.Cells.Item(15, 11).Value = -56.521739130434
$ws.Cells.Item(15, 13).Value = -47.368421052631
$ws.Cells.Item(16, 3).Value = 6
$ws.Cells.Item(16, 4).Value = 4
$ws.Cells.Item(16, 5).Value = 50
$ws.Cells.Item(16, 6).Value = 19
$ws.Cells.Item(16, 8).Value = -24
$ws.Cells.Item(16, 9).Value = 195
$ws.Cells.Item(16, 10).Value = 247
$ws.Cells.Item(16, 11).Value = -21.052631578947
$ws.Cells.Item(16, 12).Value = 6.55737704918
$ws.Cells.Item(16, 13).Value = -22
$ws.Cells.Item(16, 14).Value = -83.924154987634
$ws.Cells.Item(17, 3).Value = 10
$ws.Cells.Item(17, 5).Value = 42.857142857142
$ws.Cells.Item(17, 6).Value = 31
$ws.Cells.Item(17, 7).Value = 29
$ws.Cells.Item(17, 8).Value = 6.896551724137
$ws.Cells.Item(17, 9).Value = 299
$ws.Cells.Item(17, 10).Value = 280
$ws.Cells.Item(17, 11).Value = 6.785714285714
$ws.Cells.Item(17, 12).Value = 19.123505976095
$ws.Cells.Item(17, 13).Value = 49.5
$ws.Cells.Item(17, 14).Value = -62.151898734177
$ws.Cells.Item(18, 3).Value = 3
$ws.Cells.Item(18, 5).Value = 50
$ws.Cells.Item(18, 6).Value = 7
$ws.Cells.Item(18, 7).Value = 12
$ws.Cells.Item(18, 8).Value = -41.666666666666
$ws.Cells.Item(18, 9).Value = 125
$ws.Cells.Item(18, 10).Value = 133
$ws.Cells.Item(18, 11).Value = -6.015037593984
$ws.Cells.Item(18, 12).Value = 23.762376237623
$ws.Cells.Item(18, 13).Value = -23.312883435582
$ws.Cells.Item(18, 14).Value = -92.977528089887
$ws.Cells.Item(19, 3).Value = 17
$ws.Cells.Item(19, 4).Value = 14
$ws.Cells.Item(19, 5).Value = 21.428571428571
$ws.Cells.Item(19, 6).Value = 52
$ws.Cells.Item(19, 7).Value = 56
$ws.Cells.Item(19, 8).Value = -7.142857142857
$ws.Cells.Item(19, 9).Value = 484
$ws.Cells.Item(19, 10).Value = 524
$ws.Cells.Item(19, 11).Value = -7.633587786259
$ws.Cells.Item(19, 12).Value = -5.836575875486
$ws.Cells.Item(19, 13).Value = 47.112462006079
$ws.Cells.Item(19, 14).Value = -55.350553505535
$ws.Cells.Item(20, 3).Value = 5
$ws.Cells.Item(20, 4).Value = 2
$ws.Cells.Item(20, 5).Value = 150
$ws.Cells.Item(20, 9).Value = 203
$ws.Cells.Item(20, 10).Value = 264
$ws.Cells.Item(20, 11).Value = -23.10606060606
$ws.Cells.Item(20, 12).Value = 10.928961748633
$ws.Cells.Item(20, 13).Value = 160.25641025641
$ws.Cells.Item(20, 14).Value = -87.704421562689
$ws.Cells.Item(21, 3).Value = 41
$ws.Cells.Item(21, 5).Value = 36.666666666666
$ws.Cells.Item(21, 6).Value = 131
$ws.Cells.Item(21, 7).Value = 135
$ws.Cells.Item(21, 8).Value = -2.962962962962
$ws.Cells.Item(21, 9).Value = 1319
$ws.Cells.Item(21, 10).Value = 1475
$ws.Cells.Item(21, 11).Value = -10.57627118644
$ws.Cells.Item(21, 12).Value = 5.015923566878
$ws.Cells.Item(21, 13).Value = 26.462128475551
$ws.Cells.Item(21, 14).Value = -80.150489089541
$ws.Cells.Item(22, 6).Value = 4
$ws.Cells.Item(22, 7).Value = 2
$ws.Cells.Item(22, 8).Value = 100
$ws.Cells.Item(22, 9).Value = 37
$ws.Cells.Item(22, 11).Value = 42.307692307692
$ws.Cells.Item(22, 12).Value = 60.869565217391
$ws.Cells.Item(22, 13).Value = 68.181818181818
$ws.Cells.Item(23, 13).Value = 38.095238095238
$ws.Cells.Item(24, 3).Value = 34
$ws.Cells.Item(24, 4).Value = 23
$ws.Cells.Item(24, 5).Value = 47.826086956521
$ws.Cells.Item(24, 6).Value = 116
$ws.Cells.Item(24, 8).Value = 7.407407407407
$ws.Cells.Item(24, 9).Value = 1210
$ws.Cells.Item(24, 10).Value = 1240
$ws.Cells.Item(24, 11).Value = -2.419354838709
$ws.Cells.Item(24, 12).Value = 52.777777777777
$ws.Cells.Item(24, 13).Value = 115.302491103203
$ws.Cells.Item(25, 3).Value = 11
$ws.Cells.Item(25, 4).Value = 15
$ws.Cells.Item(25, 5).Value = -26.666666666666
$ws.Cells.Item(25, 6).Value = 46
$ws.Cells.Item(25, 7).Value = 42
$ws.Cells.Item(25, 8).Value = 9.523809523809
$ws.Cells.Item(25, 9).Value = 473
$ws.Cells.Item(25, 10).Value = 400
$ws.Cells.Item(25, 11).Value = 18.25
$ws.Cells.Item(25, 12).Value = 31.024930747922
$ws.Cells.Item(25, 13).Value = -1.663201663201
$ws.Cells.Item(26, 6).Value = 2
$ws.Cells.Item(26, 8).Value = 0
$ws.Cells.Item(26, 10).Value = 35
$ws.Cells.Item(26, 11).Value = -42.857142857142
$ws.Cells.Item(26, 12).Value = -20
$ws.Cells.Item(27, 3).Value = 3
$ws.Cells.Item(27, 5).Value = 200
$ws.Cells.Item(27, 6).Value = 5
$ws.Cells.Item(27, 7).Value = 8
$ws.Cells.Item(27, 8).Value = -37.5
$ws.Cells.Item(27, 9).Value = 49
$ws.Cells.Item(27, 10).Value = 46
$ws.Cells.Item(27, 11).Value = 6.521739130434
$ws.Cells.Item(27, 12).Value = -9.259259259259
$ws.Cells.Item(28, 8).Value = -100
$ws.Cells.Item(29, 8).Value = -100

# --- Cells changing from text placeholder to numeric value ---
$ws.Cells.Item(15, 4).Value = 1
$ws.Cells.Item(15, 4).NumberFormat = "#,##0"
$ws.Cells.Item(15, 5).Value = -100
$ws.Cells.Item(15, 5).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(26, 4).Value = 1
$ws.Cells.Item(26, 4).NumberFormat = "#,##0"
$ws.Cells.Item(26, 5).Value = -100
$ws.Cells.Item(26, 5).NumberFormat = "#,##0.0;""-""#,##0.0"

# --- Cells changing from numeric value to text placeholder ("0" or "***.*") ---
# Use a source cell that already has the desired text-placeholder style (style 14) to copy formats from,
# and force literal text entry (apostrophe prefix) only when the literal looks numeric.
$ws.Cells.Item(23, 3).Value = "'0"
$ws.Cells.Item(14, 3).Copy()
$ws.Cells.Item(23, 3).PasteSpecial(-4122)
$ws.Cells.Item(23, 4).Value = "'0"
$ws.Cells.Item(14, 3).Copy()
$ws.Cells.Item(23, 4).PasteSpecial(-4122)
$ws.Cells.Item(23, 5).Value = "***.*"
$ws.Cells.Item(14, 5).Copy()
$ws.Cells.Item(23, 5).PasteSpecial(-4122)
$ws.Cells.Item(28, 6).Value = "'0"
$ws.Cells.Item(14, 3).Copy()
$ws.Cells.Item(28, 6).PasteSpecial(-4122)
$ws.Cells.Item(29, 6).Value = "'0"
$ws.Cells.Item(14, 3).Copy()
$ws.Cells.Item(29, 6).PasteSpecial(-4122)
